$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for columns B and C (and the one-off E2 / D2 change) for rows 2-26.
# Column A already holds sequential numbers 0..6 for rows 2-8; rows 9-26 need
# new sequential numbers 7..24 filled in, with the same style as the existing
# A column cells (style index 1 -> copy from A8).

$rows = @(
    @{ Row = 2;  B = "NSE:ANGELONE";   C = "NSE:APOLLOTYRE"; D = $null; E = "NSE:ICICIPRULI"; F = $null },
    @{ Row = 3;  B = "NSE:ATALREAL";   C = "NSE:GHCL";        D = $null; E = $null;             F = $null },
    @{ Row = 4;  B = "NSE:BSE";        C = "NSE:GPIL";        D = $null; E = $null;             F = $null },
    @{ Row = 5;  B = "NSE:CAMS";       C = "NSE:GPPL";        D = $null; E = $null;             F = $null },
    @{ Row = 6;  B = "NSE:CREATIVEYE"; C = "NSE:HARDWYN";     D = $null; E = $null;             F = $null },
    @{ Row = 7;  B = "NSE:FINPIPE";    C = "NSE:HONDAPOWER";  D = $null; E = $null;             F = $null },
    @{ Row = 8;  B = "NSE:GRAVITA";    C = "NSE:IOLCP";       D = $null; E = $null;             F = $null },
    @{ Row = 9;  B = "NSE:HATHWAY";    C = "NSE:MIDHANI";     D = $null; E = $null;             F = $null },
    @{ Row = 10; B = "NSE:HIMATSEIDE"; C = "NSE:PNBGILTS";    D = $null; E = $null;             F = $null },
    @{ Row = 11; B = "NSE:HITECHCORP"; C = "NSE:POONAWALLA";  D = $null; E = $null;             F = $null },
    @{ Row = 12; B = "NSE:HLVLTD";     C = "NSE:RITES";       D = $null; E = $null;             F = $null },
    @{ Row = 13; B = "NSE:JISLJALEQS"; C = "NSE:SAMBHAAV";    D = $null; E = $null;             F = $null },
    @{ Row = 14; B = "NSE:JSWHL";      C = $null;             D = $null; E = $null;             F = $null },
    @{ Row = 15; B = "NSE:JWL";        C = $null;             D = $null; E = $null;             F = $null },
    @{ Row = 16; B = "NSE:KIRLOSENG";  C = $null;             D = $null; E = $null;             F = $null },
    @{ Row = 17; B = "NSE:MANAKSTEEL"; C = $null;             D = $null; E = $null;             F = $null },
    @{ Row = 18; B = "NSE:ORIENTLTD";  C = $null;             D = $null; E = $null;             F = $null },
    @{ Row = 19; B = "NSE:PARAGMILK";  C = $null;             D = $null; E = $null;             F = $null },
    @{ Row = 20; B = "NSE:PGIL";       C = $null;             D = $null; E = $null;             F = $null },
    @{ Row = 21; B = "NSE:POCL";       C = $null;             D = $null; E = $null;             F = $null },
    @{ Row = 22; B = "NSE:PREMIERPOL"; C = $null;             D = $null; E = $null;             F = $null },
    @{ Row = 23; B = "NSE:PRIMESECU";  C = $null;             D = $null; E = $null;             F = $null },
    @{ Row = 24; B = "NSE:RAYMOND";    C = $null;             D = $null; E = $null;             F = $null },
    @{ Row = 25; B = "NSE:RITCO";      C = $null;             D = $null; E = $null;             F = $null },
    @{ Row = 26; B = "NSE:SAGARDEEP";  C = $null;             D = $null; E = $null;             F = $null }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Column A: sequential index (row-2), same style as the existing A column.
    $ws.Cells.Item($rowNum, 1).Value = $rowNum - 2

    if ($r.B -ne $null) { $ws.Cells.Item($rowNum, 2).Value = $r.B } else { $ws.Cells.Item($rowNum, 2).Value = "" }
    if ($r.C -ne $null) { $ws.Cells.Item($rowNum, 3).Value = $r.C } else { $ws.Cells.Item($rowNum, 3).Value = "" }
    if ($r.D -ne $null) { $ws.Cells.Item($rowNum, 4).Value = $r.D } else { $ws.Cells.Item($rowNum, 4).Value = "" }
    if ($r.E -ne $null) { $ws.Cells.Item($rowNum, 5).Value = $r.E } else { $ws.Cells.Item($rowNum, 5).Value = "" }
    if ($r.F -ne $null) { $ws.Cells.Item($rowNum, 6).Value = $r.F } else { $ws.Cells.Item($rowNum, 6).Value = "" }
}

# Copy the style used on the existing A-column numeric cells (A2:A8) down to
# the newly added A9:A26 cells so they keep the same formatting.
$ws.Range("A2:A8").Copy()
$ws.Range("A9:A26").PasteSpecial(-4122)  # xlPasteFormats
